# Repro of "fixed for some error configure file" commit.
#
# The Property sheet's boolean "View" flags (columns C:F, TRUE/FALSE list
# values) were corrected: rows 68-74 gained an explicit F (View) value of
# FALSE, and the last row (75) had its Public/Private/Save/View flags reset
# to FALSE. The TRUE/FALSE list data-validation that used to cover column F
# (skipping F11/F12) now also covers the newly-edited C75:E75 block. The
# previously-active sheet/cell selection also moved back to the Property
# sheet, scrolled down near row 35, with C75 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# --- Correct the boolean flags for rows 68-75 --------------------------
$ws.Range("F68:F74").Value = $false

$ws.Range("C75:E75").Value = $false
$ws.Range("F75").Value = $false

# --- Data validation: extend the TRUE/FALSE list validation ------------
# Re-apply the existing "F" column validation (unchanged range/content)
# and add the same list validation to the newly-edited C75:E75 cells so
# they get the same TRUE/FALSE dropdown as the rest of the boolean columns.
$fRange = $ws.Range("F2:F1048576")
$fRange.Validation.Delete()
$fRange.Validation.Add(3, 1, 1, '"TRUE,FALSE"')

$newRange = $ws.Range("C75:E75")
$newRange.Validation.Delete()
$newRange.Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# --- Selection / active sheet housekeeping ------------------------------
# Work resumed on the Property sheet, scrolled so row 35 is at the top,
# with C75 (the last edited cell) selected - this also makes Property the
# workbook's active tab again (it had drifted to another sheet before).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 35
$win.ScrollColumn = 1
$ws.Range("C75").Select()

$wb.Save()
